$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "Efficiency_SOEC" worksheet between
#    "Efficiency_Alkaline" and "Sources".
# ------------------------------------------------------------------
$alkaline = $wb.Worksheets.Item("Efficiency_Alkaline")
$soec = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $alkaline)
$soec.Name = "Efficiency_SOEC"

# Headers (reuse the same text as the other efficiency sheets).
$soec.Range("A1").Value = "Power [%]"
$soec.Range("B1").Value = "Efficiency [%]"

# Power / efficiency data pairs.
$soec.Cells.Item(2, 1).Value = 0.089
$soec.Cells.Item(2, 2).Value = 3.452
$soec.Cells.Item(3, 1).Value = 0.098
$soec.Cells.Item(3, 2).Value = 16.11
$soec.Cells.Item(4, 1).Value = 0.106
$soec.Cells.Item(4, 2).Value = 26.466
$soec.Cells.Item(5, 1).Value = 0.117
$soec.Cells.Item(5, 2).Value = 37.973
$soec.Cells.Item(6, 1).Value = 0.129
$soec.Cells.Item(6, 2).Value = 47.178
$soec.Cells.Item(7, 1).Value = 0.14
$soec.Cells.Item(7, 2).Value = 54.082
$soec.Cells.Item(8, 1).Value = 0.155
$soec.Cells.Item(8, 2).Value = 60.986
$soec.Cells.Item(9, 1).Value = 0.167
$soec.Cells.Item(9, 2).Value = 65.589
$soec.Cells.Item(10, 1).Value = 0.182
$soec.Cells.Item(10, 2).Value = 70.192
$soec.Cells.Item(11, 1).Value = 0.197
$soec.Cells.Item(11, 2).Value = 73.644
$soec.Cells.Item(12, 1).Value = 0.216
$soec.Cells.Item(12, 2).Value = 77.096
$soec.Cells.Item(13, 1).Value = 0.235
$soec.Cells.Item(13, 2).Value = 79.397
$soec.Cells.Item(14, 1).Value = 0.254
$soec.Cells.Item(14, 2).Value = 80.548
$soec.Cells.Item(15, 1).Value = 0.273
$soec.Cells.Item(15, 2).Value = 81.699
$soec.Cells.Item(16, 1).Value = 0.292
$soec.Cells.Item(16, 2).Value = 82.849
$soec.Cells.Item(17, 1).Value = 0.314
$soec.Cells.Item(17, 2).Value = 82.849
$soec.Cells.Item(18, 1).Value = 0.341
$soec.Cells.Item(18, 2).Value = 82.849
$soec.Cells.Item(19, 1).Value = 0.367
$soec.Cells.Item(19, 2).Value = 82.849
$soec.Cells.Item(20, 1).Value = 0.398
$soec.Cells.Item(20, 2).Value = 82.849
$soec.Cells.Item(21, 1).Value = 0.443
$soec.Cells.Item(21, 2).Value = 82.849
$soec.Cells.Item(22, 1).Value = 0.489
$soec.Cells.Item(22, 2).Value = 81.699
$soec.Cells.Item(23, 1).Value = 0.542
$soec.Cells.Item(23, 2).Value = 81.699
$soec.Cells.Item(24, 1).Value = 0.587
$soec.Cells.Item(24, 2).Value = 80.548
$soec.Cells.Item(25, 1).Value = 0.64
$soec.Cells.Item(25, 2).Value = 79.397
$soec.Cells.Item(26, 1).Value = 0.708
$soec.Cells.Item(26, 2).Value = 78.247
$soec.Cells.Item(27, 1).Value = 0.765
$soec.Cells.Item(27, 2).Value = 77.096
$soec.Cells.Item(28, 1).Value = 0.818
$soec.Cells.Item(28, 2).Value = 75.945
$soec.Cells.Item(29, 1).Value = 0.867
$soec.Cells.Item(29, 2).Value = 74.795
$soec.Cells.Item(30, 1).Value = 0.924
$soec.Cells.Item(30, 2).Value = 73.644
$soec.Cells.Item(31, 1).Value = 0.962
$soec.Cells.Item(31, 2).Value = 72.493
$soec.Cells.Item(32, 1).Value = 0.991
$soec.Cells.Item(32, 2).Value = 71.5

# ------------------------------------------------------------------
# 2. Add the SOEC sources to the "Sources" sheet.
# ------------------------------------------------------------------
$sources = $wb.Worksheets.Item("Sources")

$sources.Range("A4").Value = "SOEC: "
$sources.Range("B4").Value = "https://ieeexplore.ieee.org/abstract/document/9025002"
$sources.Hyperlinks.Add($sources.Range("B4"), "https://ieeexplore.ieee.org/abstract/document/9025002")
# Match the existing hyperlink-cell styling (re-use the format already used by B2/B3).
$sources.Range("B3").Copy()
$sources.Range("B4").PasteSpecial(-4122)

$sources.Range("H4").Value = "and for the max efficiency: "
$sources.Range("I4").Value = "https://www.sunfire.de/files/sunfire/images/content/Produkte_Technologie/factsheets/Sunfire-Factsheet-HyLink-SOEC_2023Nov.pdf"
$sources.Columns.Item(8).ColumnWidth = 25.17

# ------------------------------------------------------------------
# 3. Update sheet selections / active sheet.
# ------------------------------------------------------------------
$pem = $wb.Worksheets.Item("Efficiency_PEM")
$pem.Activate()
$pem.Range("A1:B1").Select()

$sources.Activate()
$sources.Range("H9").Select()

$soec.Activate()
$soec.Range("E8").Select()
